$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C ("Förändrad") date values from 2023-09-16 (serial 45185)
# to 2023-10-05 (serial 45204) for rows 2 through 11.
for ($row = 2; $row -le 11; $row++) {
    $ws.Cells.Item($row, 3).Value = 45204
}
